# Apply the changes described by the diff:
# - On the "Login Details" sheet, change B5 from "secret_sauce" to "Cele"
#   and B7 from "secret_sauce" to "Nkosi".
# - Restore the view state (zoom, topLeftCell, selection) to match the
#   committed workbook.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Login Details")

$ws.Range("B5").Value = "Cele"
$ws.Range("B7").Value = "Nkosi"

# Make "Login Details" the active sheet/tab, matching activeTab="1" /
# tabSelected="1" in the target workbook.
$ws.Activate()

# Update the view: zoomed to 180%, scrolled so A607 is the top-left
# visible cell, with E612 selected.
$excel.ActiveWindow.Zoom = 180
$ws.Range("E612").Select()
$excel.ActiveWindow.ScrollRow = 607
$excel.ActiveWindow.ScrollColumn = 1

# Column A width tweak (12.5546875 -> 12.5 characters). The engine rounds
# ColumnWidth to pixel granularity and adds the default-font padding, so
# 11.64 is the input that serializes out to width="12.5".
$ws.Columns.Item(1).ColumnWidth = 11.64
